$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "firsts"
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$newSheet.Name = "Follows"
foreach ($ws in $wb.Worksheets) {
    Write-Host $ws.Name
}
